$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs)
$ws.Range("G2").Value = 0.006530999999999999
$ws.Range("H2").Value = 0.019593
$ws.Range("M2").Value = 14.25737566666667
$ws.Range("N2").Value = 42.772127
$ws.Range("O2").Value = 0.2087950866344732
$ws.Range("P2").Value = 0.2087950866344732
$ws.Range("Q2").Value = 0.093114920479
$ws.Range("R2").Value = 0.838034284311
$ws.Range("S2").Value = 0.2087950866344732
$ws.Range("T2").Value = 0.2087950866344732

# Row 3 (FAPs)
$ws.Range("G3").Value = 0.006530999999999999
$ws.Range("H3").Value = 0.019593
$ws.Range("M3").Value = 29.04277733333333
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4253229592313036
$ws.Range("P3").Value = 0.4253229592313036
$ws.Range("Q3").Value = 0.189678378764
$ws.Range("R3").Value = 1.707105408876
$ws.Range("S3").Value = 0.4253229592313036
$ws.Range("T3").Value = 0.4253229592313036

# Row 4 (MuSCs)
$ws.Range("G4").Value = 0.006530999999999999
$ws.Range("H4").Value = 0.019593
$ws.Range("M4").Value = 20.11084633333333
$ws.Range("N4").Value = 60.332539
$ws.Range("O4").Value = 0.2945174484164121
$ws.Range("P4").Value = 0.2945174484164122
$ws.Range("Q4").Value = 0.131343937403
$ws.Range("R4").Value = 1.182095436627
$ws.Range("S4").Value = 0.2945174484164121
$ws.Range("T4").Value = 0.2945174484164122

# Row 5 (Resolving-Mac)
$ws.Range("G5").Value = 0.006530999999999999
$ws.Range("H5").Value = 0.019593
$ws.Range("M5").Value = 4.873057999999999
$ws.Range("N5").Value = 14.619174
$ws.Range("O5").Value = 0.07136450571781097
$ws.Range("P5").Value = 0.07136450571781099
$ws.Range("Q5").Value = 0.031825941798
$ws.Range("R5").Value = 0.286433476182
$ws.Range("S5").Value = 0.07136450571781097
$ws.Range("T5").Value = 0.07136450571781099
